$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B,C,D,E,F,H,I,L,N across rows 2-25 (case with 380 kV)
$newData = @{
    2 = @(2.079490228807003, 0.2879520647793754, 0.06996198002805443, 0.04325744321655156, 1.816407904982455, 0.07973214163530429, 1.289826817062973, 0.2794815660180774, 1.540167290791132)
    3 = @(1.941351996942331, 0.2515369355396331, 0.07080730756321429, 0.04344196918472765, 1.780395070784337, 0.07973214163530429, 1.281869319482148, 0.2688414106235797, 1.559678060171901)
    4 = @(1.857655497802853, 0.2292199559320522, 0.07135918046975931, 0.04356372321499169, 1.759561681787176, 0.07973214163530429, 1.277813728639593, 0.2624823816003499, 1.572281759589021)
    5 = @(1.823828780032329, 0.2201354367023214, 0.07159230181238385, 0.043615468490799, 1.751391497944411, 0.07973214163530429, 1.276368864964539, 0.2599345921771175, 1.577574367419821)
    6 = @(1.818228762533977, 0.2186275314161605, 0.07163150758879766, 0.04362418947921443, 1.750054091567677, 0.07973214163530429, 1.276141470115576, 0.2595141591015278, 1.578462640997952)
    7 = @(1.857198166068258, 0.2290974000922574, 0.07136229114763992, 0.04356441244241127, 1.759450204382063, 0.07973214163530429, 1.277793402538613, 0.2624478450871948, 1.572352504460358)
    8 = @(2.031626506242048, 0.2753869384100653, 0.07024661597120385, 0.04331931523088395, 1.803724114545005, 0.07973214163530429, 1.286910138741398, 0.275776584164376, 1.546764843106903)
    9 = @(2.382666600762377, 0.3665324717064209, 0.06832062168484043, 0.04290560281572331, 1.900785110235006, 0.07973214163530429, 1.311423332185967, 0.303306419313671, 1.501556601929657)
    10 = @(2.646212965969823, 0.4337833333534036, 0.06706697319196309, 0.04264223187759963, 1.978480094564162, 0.07973214163530429, 1.333547403678651, 0.3243992698671292, 1.471394963443267)
    11 = @(2.767367854810914, 0.4644530178135824, 0.0665320856101701, 0.04253118571281789, 2.015243420357393, 0.07973214163530429, 1.344520959197794, 0.3341872422334404, 1.458341714217532)
    12 = @(2.813430589062136, 0.4760788523373662, 0.06633466289599355, 0.04249039213721018, 2.029371260761252, 0.07973214163530429, 1.348808300675643, 0.3379216943103813, 1.45349524417486)
    13 = @(2.803501940289266, 0.4735744786960936, 0.06637695276659983, 0.04249912189434069, 2.026319363222711, 0.07973214163530429, 1.347879060434593, 0.3371161657774735, 1.454534717163504)
    14 = @(2.771153764761266, 0.4654092391195945, 0.06651574054800946, 0.04252780442077331, 2.016401577414882, 0.07973214163530429, 1.344871031353108, 0.3344939160823799, 1.457941054858992)
    15 = @(2.751363594034103, 0.4604093674287242, 0.06660142093531718, 0.04254553693709218, 2.010353589485391, 0.07973214163530429, 1.343045737650556, 0.3328913609314839, 1.460040119857752)
    16 = @(2.638320852910226, 0.4317806158770736, 0.0671026448755363, 0.04264966508597157, 1.976106279781561, 0.07973214163530429, 1.332848649618413, 0.3237635018909515, 1.472261500037071)
    17 = @(2.569298471939362, 0.4142380929313276, 0.0674192230299191, 0.0427157865920873, 1.955461817716525, 0.07973214163530429, 1.326826749837608, 0.318213391306486, 1.479930244748985)
    18 = @(2.529717698968966, 0.4041552884353905, 0.06760464051460602, 0.04275464287706798, 1.943721141895651, 0.07973214163530429, 1.323448621175956, 0.3150392400651043, 1.484403893059394)
    19 = @(2.516336723212646, 0.4007426340751294, 0.06766799053170658, 0.04276794072040202, 1.939768806919034, 0.07973214163530429, 1.322319496233192, 0.3139676326939025, 1.48592936140831)
    20 = @(2.576633684742205, 0.4161047742751975, 0.06738517783020015, 0.04270866248724836, 1.95764562585768, 0.07973214163530429, 1.327458933222601, 0.3188023323657774, 1.479107391488306)
    21 = @(2.780650198641524, 0.4678072385380574, 0.0664748357553222, 0.04251934556990689, 2.019309052722775, 0.07973214163530429, 1.345750973581218, 0.3352633735926105, 1.456937907273073)
    22 = @(2.915060244947938, 0.5016674892331139, 0.06590977970333967, 0.04240294265122957, 2.060813595457461, 0.07973214163530429, 1.358475205469475, 0.3461847018055124, 1.44301166995583)
    23 = @(2.84322420634885, 0.4835889814484631, 0.06620861213769658, 0.04246439961761173, 2.038550929903124, 0.07973214163530429, 1.351613265308458, 0.3403407786027657, 1.45039268880555)
    24 = @(2.573317120232048, 0.4152608395244783, 0.06740055903872033, 0.04271188067154741, 1.956657927362158, 0.07973214163530429, 1.327172861741374, 0.3185360201726581, 1.47947920169748)
    25 = @(2.286722651238563, 0.3418292674201098, 0.06881344343145201, 0.04301038084987585, 1.873416736566369, 0.07973214163530429, 1.304074413376526, 0.2957079857949196, 1.513252251360754)
}

$cols = @("B", "C", "D", "E", "F", "H", "I", "L", "N")

foreach ($row in $newData.Keys) {
    $values = $newData[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
